$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "TC_010"
$ws.Range("C13").Value = "standard_user"
$ws.Range("D13").Value = "secret_sauce"
$ws.Range("K13").Value = "za"

$ws.Range("I17").Select()
